$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 11 - Day 8 already present in B11; add description + date
$ws.Range("C11").Value = "Working with arrays and linked lists."
$ws.Range("D11").Value = 44099

# Row 12 - Day 9 already present in B12; add description + date
$ws.Range("C12").Value = "Working with arrays and linked lists."
$ws.Range("D12").Value = 44100

# Row 13 - Day 10
$ws.Range("B13").Value = "Day 10"
$ws.Range("C13").Value = "Reading the book ""Python Crash Course""."
$ws.Range("D13").Value = 44101

# Row 14 - Day 11
$ws.Range("B14").Value = "Day 11"
$ws.Range("C14").Value = "Finish lesson about arrays and linked lists."
$ws.Range("D14").Value = 44102

# Row 15 - Day 12
$ws.Range("B15").Value = "Day 12"
$ws.Range("C15").Value = "Finish lesson about stacks and queues."
$ws.Range("D15").Value = 44103

# Row 16 - Day 13
$ws.Range("B16").Value = "Day 13"

# Row 17 - Day 14 (label entered before row 16's description, to match original authoring order)
$ws.Range("B17").Value = "Day 14"

$ws.Range("C16").Value = "Working on the lesson about recursion"
$ws.Range("C17").Value = "Working on the lesson about trees"

$ws.Range("D16").Value = 44104
$ws.Range("D17").Value = 44105

# Row 18 - Day 15
$ws.Range("B18").Value = "Day 15"
$ws.Range("C18").Value = "Working on the lesson about trees"
$ws.Range("D18").Value = 44106

# Row 19 - Day 16
$ws.Range("B19").Value = "Day 16"

# Row 20 - Day 17 (label entered before row 19's description, to match original authoring order)
$ws.Range("B20").Value = "Day 17"

$ws.Range("C19").Value = "Finish lesson about trees"
$ws.Range("C20").Value = "Working on the lesson about hashing"

$ws.Range("D19").Value = 44107
$ws.Range("D20").Value = 44108

# Apply the same date number format (style) used by existing date cells in column D,
# reusing the existing style index instead of creating a new one.
$ws.Range("D4").Copy()
$ws.Range("D11:D20").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Final selection as left by the edit
$ws.Range("C20").Select()
